# Generate Report for Archive
#
# Two changes, mirroring a localization-status refresh:
#  1) The per-locale status cells that said "Ready for handoff" now read
#     "In Translation" (Overview!E2:F2, zh-cn!C2, de-de!C2).
#  2) The "Status" column(s) that held that text are narrowed now that the
#     shorter label fits (Overview columns E/F, and column C on the two
#     per-locale detail sheets).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Compare with the literal on the left and an explicit [string] cast
        # on the right: some cells hold text that looks boolean ("True"/
        # "False") and this host's Value2/Text getters type-convert those to
        # real Booleans, which -eq would otherwise silently coerce the
        # string literal into (any non-empty string -> $true) and cause a
        # false-positive match.
        if ($oldStatus -eq [string]$cell.Text) {
            $cell.Value = $newStatus
        }
    }
}

# Narrow the status columns (closest COM-reachable width to 13.4101845877511
# on this engine's 1/6-pt ColumnWidth grid).
$newColWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newColWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newColWidth
